$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "313/422"
$ws.Range("H3").Value = "412/422"
$ws.Range("H8").Value = "144/422"
$ws.Range("H12").Value = "315/422"
$ws.Range("H14").Value = "57/422"
$ws.Range("H20").Value = "270/422"
$ws.Range("H28").Value = "386/422"
$ws.Range("H38").Value = "307/422"
$ws.Range("H41").Value = "396/422"
$ws.Range("H42").Value = "224/422"
$ws.Range("H45").Value = "308/422"
$ws.Range("H47").Value = "406/422"
$ws.Range("H50").Value = "302/422"
$ws.Range("H56").Value = "297/422"
$ws.Range("H65").Value = "173/422"
$ws.Range("H73").Value = "344/422"
$ws.Range("H75").Value = "174/422"
$ws.Range("H76").Value = "286/422"
$ws.Range("H85").Value = "250/422"
$ws.Range("H87").Value = "388/422"
$ws.Range("H89").Value = "340/422"
$ws.Range("H90").Value = "365/422"
$ws.Range("H92").Value = "368/422"
$ws.Range("H93").Value = "379/422"
